$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = -4
$ws.Range("F5").Value = -2
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = 0
